$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 187.0867421867131
$ws.Range("R2").Value = 1683.780679680418
$ws.Range("S2").Value = 0.002476505008608346
$ws.Range("T2").Value = 0.002476505008608346

$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 87.17920859167909
$ws.Range("R3").Value = 784.6128773251119
$ws.Range("S3").Value = 0.001154008799342588
$ws.Range("T3").Value = 0.001154008799342588

$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 24.42531779601222
$ws.Range("R4").Value = 219.82786016411
$ws.Range("S4").Value = 0.0003233228669849103
$ws.Range("T4").Value = 0.0003233228669849103

$ws.Range("G5").Value = 1.468507333333333
$ws.Range("H5").Value = 4.405521999999999
$ws.Range("I5").Value = 0.005118279455112885
$ws.Range("J5").Value = 0.005118279455112885
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 87.96744018239711
$ws.Range("R5").Value = 791.706961641574
$ws.Range("S5").Value = 0.001164442780177041
$ws.Range("T5").Value = 0.001164442780177041

$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 33066.54993862413
$ws.Range("R6").Value = 297598.9494476172
$ws.Range("S6").Value = 0.4377086028825855
$ws.Range("T6").Value = 0.4377086028825856

$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("Q7").Value = 15408.44434411888
$ws.Range("R7").Value = 138675.9990970699
$ws.Range("S7").Value = 0.203964691175126
$ws.Range("T7").Value = 0.203964691175126

$ws.Range("I8").Value = 0.9046276674881553
$ws.Range("J8").Value = 0.9046276674881553
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 4317.040220105789
$ws.Range("R8").Value = 38853.36198095209
$ws.Range("S8").Value = 0.05714553368397341
$ws.Range("T8").Value = 0.05714553368397341

$ws.Range("I9").Value = 0.9046276674881553
$ws.Range("J9").Value = 0.9046276674881553
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 15547.75993085173
$ws.Range("R9").Value = 139929.8393776655
$ws.Range("S9").Value = 0.2058088397464704
$ws.Range("T9").Value = 0.2058088397464704

$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 75.04750459246443
$ws.Range("R10").Value = 675.42754133218
$ws.Range("S10").Value = 0.0009934189822029809
$ws.Range("T10").Value = 0.0009934189822029813

$ws.Range("G11").Value = 0.5890733333333333
$ws.Range("H11").Value = 1.76722
$ws.Range("I11").Value = 0.002053133730501083
$ws.Range("J11").Value = 0.002053133730501083
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 34.97084817812444
$ws.Range("R11").Value = 314.7376336031199
$ws.Range("S11").Value = 0.0004629161834566273
$ws.Range("T11").Value = 0.0004629161834566275

$ws.Range("G12").Value = 0.5890733333333333
$ws.Range("H12").Value = 1.76722
$ws.Range("I12").Value = 0.002053133730501083
$ws.Range("J12").Value = 0.002053133730501083
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 9.79791046678889
$ws.Range("R12").Value = 88.1811942011
$ws.Range("S12").Value = 0.0001296969205903575
$ws.Range("T12").Value = 0.0001296969205903575

$ws.Range("G13").Value = 0.5890733333333333
$ws.Range("H13").Value = 1.76722
$ws.Range("I13").Value = 0.002053133730501083
$ws.Range("J13").Value = 0.002053133730501083
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 35.28703741330445
$ws.Range("R13").Value = 317.58333671974
$ws.Range("S13").Value = 0.0004671016442511173
$ws.Range("T13").Value = 0.0004671016442511174

$ws.Range("G14").Value = 25.306101
$ws.Range("H14").Value = 75.918303
$ws.Range("I14").Value = 0.0882009193262308
$ws.Range("J14").Value = 0.0882009193262308
$ws.Range("M14").Value = 127.3992563333333
$ws.Range("N14").Value = 382.197769
$ws.Range("O14").Value = 0.4838549810199306
$ws.Range("P14").Value = 0.4838549810199307
$ws.Range("Q14").Value = 3223.978448096223
$ws.Range("R14").Value = 29015.806032866
$ws.Range("S14").Value = 0.04267645414653384
$ws.Range("T14").Value = 0.04267645414653384

$ws.Range("G15").Value = 25.306101
$ws.Range("H15").Value = 75.918303
$ws.Range("I15").Value = 0.0882009193262308
$ws.Range("J15").Value = 0.0882009193262308
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2254681108101269
$ws.Range("P15").Value = 0.2254681108101269
$ws.Range("Q15").Value = 1502.318584077732
$ws.Range("R15").Value = 13520.86725669959
$ws.Range("S15").Value = 0.01988649465220167
$ws.Range("T15").Value = 0.01988649465220167

$ws.Range("G16").Value = 25.306101
$ws.Range("H16").Value = 75.918303
$ws.Range("I16").Value = 0.0882009193262308
$ws.Range("J16").Value = 0.0882009193262308
$ws.Range("M16").Value = 16.63275166666667
$ws.Range("N16").Value = 49.898255
$ws.Range("O16").Value = 0.06317022542837675
$ws.Range("P16").Value = 0.06317022542837675
$ws.Range("Q16").Value = 420.910093584585
$ws.Range("R16").Value = 3788.190842261265
$ws.Range("S16").Value = 0.005571671956828071
$ws.Range("T16").Value = 0.005571671956828071

$ws.Range("G17").Value = 25.306101
$ws.Range("H17").Value = 75.918303
$ws.Range("I17").Value = 0.0882009193262308
$ws.Range("J17").Value = 0.0882009193262308
$ws.Range("M17").Value = 59.90262233333334
$ws.Range("N17").Value = 179.707867
$ws.Range("O17").Value = 0.2275066827415657
$ws.Range("P17").Value = 0.2275066827415658
$ws.Range("Q17").Value = 1515.901810932189
$ws.Range("R17").Value = 13643.1162983897
$ws.Range("S17").Value = 0.02006629857066723
$ws.Range("T17").Value = 0.02006629857066723

